$wb = $excel.ActiveWorkbook

# Sheet 1 = "nhap-linhkien" (import), Sheet 2 = "xuat-linhkien" (export)
$wsNhap = $wb.Worksheets.Item(1)
$wsXuat = $wb.Worksheets.Item(2)

# --- nhap-linhkien: write the new stock-take row (row 2) -------------------
# Plain string columns first (A-E), in left-to-right order, so the shared
# string table grows in the same order as the source edit.
$wsNhap.Cells.Item(2, 1).Value = "PCB -RF_SOC 2021.2.19( PHIÊN BẢN CŨ)"
$wsNhap.Cells.Item(2, 2).Value = "PCB RF SOC 1 PHA"

# Sổ Hợp Đồng is genuinely blank text (not "no cell") in the target sheet,
# so force a literal empty string instead of letting an empty value clear
# the cell outright.
$wsNhap.Cells.Item(2, 3).Value = "'"
$wsNhap.Cells.Item(2, 3).ClearFormats()

$wsNhap.Cells.Item(2, 4).Value = "MODULE RF 1P SOC"
$wsNhap.Cells.Item(2, 5).Value = "TỒN 30/06/2021"

# Ngày Nhập looks like a date but must stay literal text, not an Excel
# date serial - force text with a leading quote, then drop the resulting
# quote-prefix style so the cell is plain text again.
$wsNhap.Cells.Item(2, 6).Value = "'2021-06-30"
$wsNhap.Cells.Item(2, 6).ClearFormats()

$wsNhap.Cells.Item(2, 7).Value = "Cái"

# Numeric columns.
$wsNhap.Cells.Item(2, 8).Value = 92
$wsNhap.Cells.Item(2, 9).Value = 0
$wsNhap.Cells.Item(2, 10).Value = 0

# --- xuat-linhkien: drop the old sample row, keep only the header ----------
$wsXuat.Range("A2:J2").Clear()
